$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A (NC)
$ws.Cells.Item(2, 1).Value = 19330051920303
$ws.Cells.Item(3, 1).Value = 18330051920253

# Column B (Paterno)
$ws.Cells.Item(2, 2).Value = "VAZQUEZ"
$ws.Cells.Item(3, 2).Value = "CONTRERAS"

# Column C (Materno)
$ws.Cells.Item(2, 3).Value = "VERA"
$ws.Cells.Item(3, 3).Value = "CASTRO"

# Column D (Nombres)
$ws.Cells.Item(2, 4).Value = "MARIA FERNANDA"
$ws.Cells.Item(3, 4).Value = "SAMUEL"

# Column E (Nombre_Largo)
$ws.Cells.Item(2, 5).Value = "REALIZA ANÁLISIS HEMATOLÓGICOS DE SERIE ROJA"
$ws.Cells.Item(3, 5).Value = "ANALIZA SANGRE MEDIANTE PRUEBAS HORMONALES, TOXICOLÓGICAS Y DE MARCADORES TUMORALES"

# Column F (Grupo/Mat)
$ws.Cells.Item(2, 6).Value = "4ALCM"
$ws.Cells.Item(3, 6).Value = "6ALCM"

# Column G (Reprobadas)
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(3, 7).Value = 2
